# geopoints.xlsx — preliminary check-in
#
# Adds a "properties" worksheet (survey -> properties -> choices -> settings)
# that will be used to generate properties.csv, alongside the existing
# survey/choices/settings sheets used for definitions.csv. The new sheet is
# inserted right after "survey" and becomes the active sheet/tab, with the
# "survey" sheet's selection moved off its old edit point.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Insert the new sheet right after "survey" (so tab order becomes
# survey, properties, choices, settings).
$props = $wb.Worksheets.Add($null, $survey)
$props.Name = "properties"

# --- Header row -------------------------------------------------------
$props.Range("A1").Value = "partition"
$props.Range("B1").Value = "aspect"
$props.Range("C1").Value = "key"
$props.Range("D1").Value = "type"
$props.Range("E1").Value = "value"

# --- Row 2: Table / colOrder ------------------------------------------
$props.Range("A2").Value = "Table"
$props.Range("B2").Value = "default"
$props.Range("C2").Value = "colOrder"
$props.Range("D2").Value = "array"
$props.Range("E2").Value = '["client_id","step","transportation_mode","transportation_mode_other","description","coordinates_latitude","coordinates_longitude","coordinates_altitude","coordinates_accuracy"]'

# --- Row 3: Table / mapListViewFileName --------------------------------
$props.Range("A3").Value = "Table"
$props.Range("B3").Value = "default"
$props.Range("C3").Value = "mapListViewFileName"
$props.Range("D3").Value = "configpath"
$props.Range("E3").Value = "config/tables/geopoints/html/geopoints_map_list.html"

# --- Row 4: Table / defaultViewType ------------------------------------
$props.Range("A4").Value = "Table"
$props.Range("B4").Value = "default"
$props.Range("C4").Value = "defaultViewType"
$props.Range("D4").Value = "string"
$props.Range("E4").Value = "MAP"

# --- Rows 5-7: TableMapFragment (filled column-by-column, matching the
#     original authoring order so shared-string interning order lines up)
$props.Range("A5").Value = "TableMapFragment"
$props.Range("A6").Value = "TableMapFragment"
$props.Range("A7").Value = "TableMapFragment"

$props.Range("B5").Value = "default"
$props.Range("B6").Value = "default"
$props.Range("B7").Value = "default"

$props.Range("C5").Value = "keyColorRuleType"
$props.Range("C6").Value = "keyMapLatCol"
$props.Range("C7").Value = "keyMapLongCol"

$props.Range("D5").Value = "string"
$props.Range("D6").Value = "string"
$props.Range("D7").Value = "string"

$props.Range("E5").Value = "None"
$props.Range("E6").Value = "coordinates_latitude"
$props.Range("E7").Value = "coordinates_longitude"

# --- Selections ---------------------------------------------------------
# "survey" is no longer the tab in focus; its selection moved too.
$null = $survey.Range("E29").Select()

# "properties" becomes the active/selected tab.
$null = $props.Activate()
$null = $props.Range("B9").Select()
